$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2022/3/20完成" progress note for the "move several pages in
# one frame" optimization task (row 8), next to the task description in A8.
$ws.Range("B8").Value = "2022/3/20完成"

# Update the sheet's active selection to the newly filled cell.
$ws.Range("B8").Select()

# Recolor the theme's "Background 1" / Light 1 color (used as the sheet's
# background tint) from white to a light green.
$tcs = $wb.Theme.ThemeColorScheme
$rgbVal = 0xCC + (0xE8 * 256) + (0xCF * 65536)
$tcs.Colors(2).RGB = $rgbVal
